# texts.xlsx - "Translation" sheet, Table8 (Text ID / Typography Name / Alignment / GB / ...)
#
# The duplicate "Brake" entry (SingleUseId32, centered) that used to live on
# row 24 is being replaced. All rows below it move up by one, and a new
# left-aligned "Brake" entry (SingleUseId35) is appended as the last row of
# the table (row 40) - this is where the new "Mute" audio-settings text will
# eventually be added to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")
$lo = $ws.ListObjects.Item("Table8")

$tableRange = $lo.Range

# Remove the old row 24 entirely; Excel shifts every row below it up by one.
$ws.Rows.Item(24).Delete()

# The table shrinks by one row when a row inside it is removed - put it back
# to its original size so the (now blank) row 40 is part of the table again.
$lo.Resize($tableRange)

# Populate the newly freed last row (row 40) with the new entry.
$ws.Range("B40").Value = "SingleUseId35"
$ws.Range("C40").Value = "SansSerif40px"
$ws.Range("D40").Value = "Left"
$ws.Range("E40").Value = "Brake"
$ws.Range("F40").Value = "LTR"
